# Append the 2025-04-05 price row (carried forward from 2025-04-04 / row 34)
# to every price sheet in the workbook, per the "Updated Argent prices" commit.

$wb = $excel.ActiveWorkbook

# Sheet name -> new Price value for the appended 2025-04-05 row.
# (Every sheet simply repeats its last known price, same as the prior
# day-over-day "carry forward" rows already in the sheet.)
$updates = [ordered]@{
    "N-Dense"                   = "40"
    "N-Type"                    = "42"
    "N-type Wafer"               = "1.28"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,487"
    "Silver Busbar front-side"   = "8,215"
    "Silver finger front-side"   = "8,265"
    "USD_CNY"                    = "7.3068"
}

$newDate = "2025-04-05"
$newRow = 35

foreach ($name in $updates.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $price = $updates[$name]

    $dateCell = $ws.Cells.Item($newRow, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate
    $dateCell.Style = "Normal"

    $priceCell = $ws.Cells.Item($newRow, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
